$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting existing rows 117-137 down to 118-138.
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new data record.
$ws.Range("A117").Value() = 6
$ws.Range("B117").Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C117").Value() = "Metropolitana"
$ws.Range("D117").Value() = 44694
$ws.Range("E117").Value() = 13
$ws.Range("F117").Value() = "Fruta"
$ws.Range("G117").Value() = 100107
$ws.Range("H117").Value() = "Otros"
$ws.Range("I117").Value() = 100107001
$ws.Range("J117").Value() = "Caqui"
$ws.Range("K117").Value() = "Fuyu"
$ws.Range("L117").Value() = "Tercera"
$ws.Range("M117").Value() = 10
$ws.Range("N117").Value() = 200000
$ws.Range("O117").Value() = 200000
$ws.Range("P117").Value() = 200000
$ws.Range("Q117").Value() = "$/bins (400 kilos)"
$ws.Range("R117").Value() = "Región de O'Higgins"
$ws.Range("S117").Value() = 500
$ws.Range("T117").Value() = 400
